$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Row value/percentage updates (no coin identity change) ---
Set-TextValue 2 4 "63.088.49"
Set-TextValue 2 5 "  -1.99%  "
Set-TextValue 3 4 "3.414.38"
Set-TextValue 3 5 "  -2.57%  "
Set-TextValue 5 4 "576.74"
Set-TextValue 5 5 "  -2.23%  "
Set-TextValue 6 4 "127.01"
Set-TextValue 6 5 "  -5.47%  "
Set-TextValue 7 5 "  +0.03%  "
Set-TextValue 8 4 "3.414.18"
Set-TextValue 8 5 "  -2.54%  "
Set-TextValue 9 4 "0.477"
Set-TextValue 9 5 "  -2.20%  "
Set-TextValue 10 4 "7.45"
Set-TextValue 10 5 "  -0.46%  "
Set-TextValue 11 4 "0.121"
Set-TextValue 11 5 "  -2.46%  "
Set-TextValue 12 4 "0.379"
Set-TextValue 12 5 "  -1.70%  "
Set-TextValue 13 4 "4.004.35"
Set-TextValue 13 5 "  -2.41%  "
Set-TextValue 14 5 "  -0.79%  "
Set-TextValue 15 4 "3.420.43"
Set-TextValue 15 5 "  -2.40%  "
Set-TextValue 16 4 "0.0000174"
Set-TextValue 16 5 "  -4.02%  "
Set-TextValue 17 4 "63.120.23"
Set-TextValue 17 5 "  -1.94%  "
Set-TextValue 18 4 "24.83"
Set-TextValue 18 5 "  -3.60%  "
Set-TextValue 19 4 "9.65"
Set-TextValue 19 5 "  -2.56%  "
Set-TextValue 20 4 "5.69"
Set-TextValue 20 5 "  -1.10%  "
Set-TextValue 21 4 "13.16"
Set-TextValue 21 5 "  -2.96%  "
Set-TextValue 22 4 "378.38"
Set-TextValue 22 5 "  -3.76%  "
Set-TextValue 23 4 "0.560"
Set-TextValue 23 5 "  -2.61%  "
Set-TextValue 24 4 "3.553.63"
Set-TextValue 24 5 "  -2.51%  "
Set-TextValue 25 4 "72.79"
Set-TextValue 25 5 "  -2.46%  "
Set-TextValue 26 5 "  -0.09%  "
Set-TextValue 27 5 "  -7.01%  "
Set-TextValue 28 5 "  -0.12%  "
Set-TextValue 29 4 "6.99"
Set-TextValue 29 5 "  -5.28%  "
Set-TextValue 30 4 "2.16"
Set-TextValue 30 5 "  -4.37%  "
Set-TextValue 31 4 "7.90"
Set-TextValue 31 5 "  -4.38%  "
Set-TextValue 32 4 "0.153"
Set-TextValue 32 5 "  -2.94%  "
Set-TextValue 33 4 "1.40"
Set-TextValue 33 5 "  -4.61%  "
Set-TextValue 34 4 "3.446.38"
Set-TextValue 34 5 "  -2.36%  "
Set-TextValue 35 5 "  -0.02%  "
Set-TextValue 36 4 "22.80"
Set-TextValue 36 5 "  -2.51%  "
Set-TextValue 37 4 "5.28"
Set-TextValue 37 5 "  -1.32%  "
Set-TextValue 38 4 "6.75"
Set-TextValue 38 5 "  -2.85%  "
Set-TextValue 39 4 "164.33"
Set-TextValue 39 5 "  -1.99%  "
Set-TextValue 40 4 "1.50"
Set-TextValue 40 5 "  -3.81%  "
Set-TextValue 41 4 "0.0761"
Set-TextValue 41 5 "  -3.41%  "
Set-TextValue 44 4 "41.79"
Set-TextValue 44 5 "  -0.90%  "
Set-TextValue 45 4 "4.28"
Set-TextValue 45 5 "  -3.59%  "
Set-TextValue 46 4 "1.59"
Set-TextValue 46 5 "  -4.92%  "
Set-TextValue 47 5 "  -8.83%  "
Set-TextValue 48 4 "1.09"
Set-TextValue 48 5 "  -7.20%  "
Set-TextValue 49 4 "6.70"
Set-TextValue 49 5 "  -1.41%  "

# --- Row swaps: coin identity + link + price + volume changed ---
Set-TextValue 42 2 "FirstDigitalUSD"
Set-TextValue 42 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue 42 4 "1.00"
Set-TextValue 42 5 "  +0.09%  "
Set-TextValue 43 2 "Mantle"
Set-TextValue 43 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 43 4 "0.781"
Set-TextValue 43 5 "  -3.57%  "
Set-TextValue 50 2 "Maker"
Set-TextValue 50 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 50 4 "2.262.57"
Set-TextValue 50 5 "  -5.21%  "
Set-TextValue 51 2 "SuiNetwork"
Set-TextValue 51 3 "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue 51 4 "0.863"
Set-TextValue 51 5 "  -3.94%  "
